$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.550.98'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.444.14'
$ws.Range("E3").Value = '  -3.44%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.34'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.00'
$ws.Range("E6").Value = '  -7.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.445.13'
$ws.Range("E7").Value = '  -3.37%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.39'
$ws.Range("E10").Value = '  -6.42%  '
$ws.Range("E12").Value = '  -7.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.026.94'
$ws.Range("E13").Value = '  -3.56%  '
$ws.Range("E14").Value = '  -10.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.66'
$ws.Range("E15").Value = '  -9.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.432.88'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.535.86'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("E19").Value = '  -9.57%  '
$ws.Range("E20").Value = '  -7.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.67'
$ws.Range("E21").Value = '  -7.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.00'
$ws.Range("E22").Value = '  -6.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.33'
$ws.Range("E24").Value = '  -5.70%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.587.09'
$ws.Range("E26").Value = '  -3.22%  '
$ws.Range("E27").Value = '  -10.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  -8.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("E30").Value = '  -9.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.22'
$ws.Range("E31").Value = '  -11.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.450.58'
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.148'
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("E35").Value = '  -7.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '171.71'
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.22'
$ws.Range("E37").Value = '  -12.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("E38").Value = '  -10.14%  '
$ws.Range("E39").Value = '  -7.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.82'
$ws.Range("E40").Value = '  -11.32%  '
$ws.Range("E41").Value = '  -8.15%  '
$ws.Range("E42").Value = '  -5.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.48'
$ws.Range("E43").Value = '  -4.93%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.44'
$ws.Range("E45").Value = '  -14.13%  '
$ws.Range("E46").Value = '  -11.96%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.83'
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.59'
$ws.Range("E50").Value = '  -15.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.204.93'
$ws.Range("E51").Value = '  -7.45%  '
